$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 (8a076200-...md) and Row 4 (9d99b3c0-...md) both share the
# "Latest HO Xliff Generate Date" value that changed from
# 2016-08-27 02:16:09 -> 2016-08-27 02:16:54
$wsOverview.Range("G3").Value = "2016-08-27 02:16:54"
$wsOverview.Range("G4").Value = "2016-08-27 02:16:54"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E) changed from "ht" to "mt" for rows 3 and 4
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H) changed from
# 2016-08-27 02:16:01 -> 2016-08-27 02:16:49
$wsZhCn.Range("H3").Value = "2016-08-27 02:16:49"
$wsZhCn.Range("H4").Value = "2016-08-27 02:16:49"
# Correspond Handback DateTime column (K) changed from
# 2016-08-27 02:16:24 -> 2016-08-27 02:17:12
$wsZhCn.Range("K3").Value = "2016-08-27 02:17:12"
$wsZhCn.Range("K4").Value = "2016-08-27 02:17:12"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# "Latest HO Xliff Generate Date" / Correspond Handoff Datetime column (H)
# matches the Overview sheet change: 2016-08-27 02:16:09 -> 2016-08-27 02:16:54
$wsDeDe.Range("H3").Value = "2016-08-27 02:16:54"
$wsDeDe.Range("H4").Value = "2016-08-27 02:16:54"
# Priority column (E) changed from "ht" to "mt" for rows 3 and 4
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handback DateTime column (K) changed from
# 2016-08-27 02:16:30 -> 2016-08-27 02:17:19
$wsDeDe.Range("K3").Value = "2016-08-27 02:17:19"
$wsDeDe.Range("K4").Value = "2016-08-27 02:17:19"
